$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.349091053009033
$ws.Range("B1").Value = 2.500438213348389
$ws.Range("C1").Value = 4.272638320922852
$ws.Range("D1").Value = 4.333022117614746
$ws.Range("E1").Value = 1.624800086021423
